$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Typography": bump two existing font sizes, then add four new
# typography rows (Typography_09 .. Typography_12) for the stopwatch work.
# ---------------------------------------------------------------------------
$typo = $wb.Worksheets.Item("Typography")

$typo.Range("D11").Value = 90
$typo.Range("D12").Value = 19

$typo.Rows("16:19").Insert()

$typo.Range("B16").Value = "Typography_09"
$typo.Range("C16").Value = "malgunbd.ttf"
$typo.Range("D16").Value = 22
$typo.Range("E16").Value = 4
$typo.Range("F16").Value = "?"

$typo.Range("B17").Value = "Typography_10"
$typo.Range("C17").Value = "malgunbd.ttf"
$typo.Range("D17").Value = 18
$typo.Range("E17").Value = 4
$typo.Range("F17").Value = "?"
$typo.Range("H17").Value = "0123456789 :APM"

$typo.Range("B18").Value = "Typography_11"
$typo.Range("C18").Value = "malgunbd.ttf"
$typo.Range("D18").Value = 13
$typo.Range("E18").Value = 4
$typo.Range("F18").Value = "?"

$typo.Range("B19").Value = "Typography_12"
$typo.Range("C19").Value = "malgunbd.ttf"
$typo.Range("D19").Value = 24
$typo.Range("E19").Value = 4
$typo.Range("F19").Value = "?"

# ---------------------------------------------------------------------------
# Sheet "Translation": drop the now-unused "(braking)"/"(gas)" rows and the
# stray "Current Time" row, retarget a couple of entries onto the new
# typographies, and append the new stopwatch-related text rows.
# ---------------------------------------------------------------------------
$tr = $wb.Worksheets.Item("Translation")

# Remove "(braking)" (row 15) and "(gas)" (row 16)
$tr.Rows("15:16").Delete()

# Remove the lone "Current Time" row (now shifted up to row 19)
$tr.Rows("19:19").Delete()

# Column F on this sheet is always free text - force text formatting so
# numeric-looking strings ("100", "0", etc.) are not coerced into numbers.
$tr.Columns("F").NumberFormat = "@"

$tr.Range("F11").Value = "100"
$tr.Range("F13").Value = "Battery"
$tr.Range("C14").Value = "Typography_09"
$tr.Range("F14").Value = "6500"
$tr.Range("C18").Value = "Typography_10"

# Append the new stopwatch text rows
$tr.Rows("23:28").Insert()

$tr.Range("B23").Value = "SingleUseId28"
$tr.Range("C23").Value = "Default"
$tr.Range("D23").Value = "Left"
$tr.Range("E23").Value = "LTR"
$tr.Range("F23").Value = "100%"

$tr.Range("B24").Value = "SingleUseId29"
$tr.Range("C24").Value = "Typography_00"
$tr.Range("D24").Value = "Center"
$tr.Range("E24").Value = "LTR"
$tr.Range("F24").Value = "rpm"

$tr.Range("B25").Value = "SingleUseId30"
$tr.Range("C25").Value = "Typography_10"
$tr.Range("D25").Value = "Center"
$tr.Range("E25").Value = "LTR"
$tr.Range("F25").Value = "Current Time"

$tr.Range("B26").Value = "SingleUseId31"
$tr.Range("C26").Value = "Typography_11"
$tr.Range("D26").Value = "Left"
$tr.Range("E26").Value = "LTR"
$tr.Range("F26").Value = "Start/Stop"

$tr.Range("B27").Value = "SingleUseId32"
$tr.Range("C27").Value = "Typography_11"
$tr.Range("D27").Value = "Left"
$tr.Range("E27").Value = "LTR"
$tr.Range("F27").Value = "Reset"

$tr.Range("B28").Value = "SingleUseId33"
$tr.Range("C28").Value = "Typography_12"
$tr.Range("D28").Value = "Left"
$tr.Range("E28").Value = "LTR"
$tr.Range("F28").Value = "99:99.99"
